# "smooth selection frame + dialog panel resize"
#
# - Add two new header columns (D: "english", E: "Croatian")
# - Resize dialogue rows 2-4 (taller/adjusted panel heights)
# - Move the active cell / selection frame from C8 to E2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New language header columns next to the existing ID / DATA / TEXT headers.
$ws.Range("D1").Value = "english"
$ws.Range("E1").Value = "Croatian"

# Resize the dialogue rows so the panel fits its (growing) contents.
$ws.Rows.Item(2).RowHeight = 122.8
$ws.Rows.Item(3).RowHeight = 33.25
$ws.Rows.Item(4).RowHeight = 29.85

# Move / smooth the selection frame onto the newly added column.
$ws.Range("E2").Select()
